# Updated cryptos list on Thu May 25 09:44:39 UTC 2023 with GitHub Actions
#
# Refreshes the Coin/Link/Price/Volume(1h) table (Sheet1, columns B:E) with
# the latest scrape. Most rows only get new Price (D) / Volume(1h) (E)
# figures, but a few coins changed rank and swapped places with their
# neighbour, so those rows also get new Coin (B) / Link (C) values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# (cell, new value) pairs taken from the refreshed scrape. 'Numeric' flags
# Price cells whose text (e.g. "1.004", "45.82") would otherwise be
# auto-coerced to a real number by a plain Value assignment, silently
# changing the cell's type/formatting - those need to be forced to text.
$updates = @(
    @{ Cell = 'D2'; Value = '26.311.06'; Numeric = $false },
    @{ Cell = 'E2'; Value = '  -2.43%  '; Numeric = $false },
    @{ Cell = 'D3'; Value = '1.791.21'; Numeric = $false },
    @{ Cell = 'E3'; Value = '  -2.14%  '; Numeric = $false },
    @{ Cell = 'D4'; Value = '1.004'; Numeric = $true },
    @{ Cell = 'E4'; Value = '  -0.25%  '; Numeric = $false },
    @{ Cell = 'E5'; Value = '  -0.40%  '; Numeric = $false },
    @{ Cell = 'D6'; Value = '306.30'; Numeric = $true },
    @{ Cell = 'E6'; Value = '  -1.66%  '; Numeric = $false },
    @{ Cell = 'D7'; Value = '0.4502'; Numeric = $true },
    @{ Cell = 'E7'; Value = '  -1.71%  '; Numeric = $false },
    @{ Cell = 'D8'; Value = '0.3591'; Numeric = $true },
    @{ Cell = 'E8'; Value = '  -3.21%  '; Numeric = $false },
    @{ Cell = 'D9'; Value = '45.82'; Numeric = $true },
    @{ Cell = 'E9'; Value = '  -0.13%  '; Numeric = $false },
    @{ Cell = 'D10'; Value = '0.07069'; Numeric = $true },
    @{ Cell = 'E10'; Value = '  -1.64%  '; Numeric = $false },
    @{ Cell = 'D11'; Value = '0.8821'; Numeric = $true },
    @{ Cell = 'E11'; Value = '  +0.56%  '; Numeric = $false },
    @{ Cell = 'D12'; Value = '0.07740'; Numeric = $true },
    @{ Cell = 'E12'; Value = '  -0.50%  '; Numeric = $false },
    @{ Cell = 'D13'; Value = '19.42'; Numeric = $true },
    @{ Cell = 'E13'; Value = '  -1.26%  '; Numeric = $false },
    @{ Cell = 'D14'; Value = '1.800.46'; Numeric = $false },
    @{ Cell = 'E14'; Value = '  -1.79%  '; Numeric = $false },
    @{ Cell = 'D15'; Value = '5.270'; Numeric = $true },
    @{ Cell = 'E15'; Value = '  -1.23%  '; Numeric = $false },
    @{ Cell = 'D16'; Value = '6.315'; Numeric = $true },
    @{ Cell = 'E16'; Value = '  -1.42%  '; Numeric = $false },
    @{ Cell = 'D17'; Value = '84.75'; Numeric = $true },
    @{ Cell = 'E17'; Value = '  -2.99%  '; Numeric = $false },
    @{ Cell = 'E18'; Value = '  -0.30%  '; Numeric = $false },
    @{ Cell = 'D19'; Value = '0.000008501'; Numeric = $true },
    @{ Cell = 'E19'; Value = '  -2.51%  '; Numeric = $false },
    @{ Cell = 'D20'; Value = '1.003'; Numeric = $true },
    @{ Cell = 'E20'; Value = '  -0.34%  '; Numeric = $false },
    @{ Cell = 'E21'; Value = '  -1.95%  '; Numeric = $false },
    @{ Cell = 'D22'; Value = '26.348.79'; Numeric = $false },
    @{ Cell = 'E22'; Value = '  -2.39%  '; Numeric = $false },
    @{ Cell = 'D23'; Value = '4.965'; Numeric = $true },
    @{ Cell = 'E23'; Value = '  -0.91%  '; Numeric = $false },
    @{ Cell = 'B24'; Value = 'Cosmos'; Numeric = $false },
    @{ Cell = 'C24'; Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'; Numeric = $false },
    @{ Cell = 'D24'; Value = '10.51'; Numeric = $true },
    @{ Cell = 'E24'; Value = '  +0.75%  '; Numeric = $false },
    @{ Cell = 'B25'; Value = 'WrappedliquidstakedEther2.0'; Numeric = $false },
    @{ Cell = 'C25'; Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'; Numeric = $false },
    @{ Cell = 'D25'; Value = '2.015.89'; Numeric = $false },
    @{ Cell = 'E25'; Value = '  -2.17%  '; Numeric = $false },
    @{ Cell = 'D26'; Value = '1.970'; Numeric = $true },
    @{ Cell = 'E26'; Value = '  -2.34%  '; Numeric = $false },
    @{ Cell = 'D27'; Value = '151.04'; Numeric = $true },
    @{ Cell = 'E27'; Value = '  -0.34%  '; Numeric = $false },
    @{ Cell = 'D28'; Value = '17.79'; Numeric = $true },
    @{ Cell = 'E28'; Value = '  -2.36%  '; Numeric = $false },
    @{ Cell = 'D29'; Value = '2.008'; Numeric = $true },
    @{ Cell = 'E29'; Value = '  +2.19%  '; Numeric = $false },
    @{ Cell = 'D30'; Value = '111.83'; Numeric = $true },
    @{ Cell = 'E30'; Value = '  -2.02%  '; Numeric = $false },
    @{ Cell = 'D31'; Value = '4.885'; Numeric = $true },
    @{ Cell = 'E31'; Value = '  -1.14%  '; Numeric = $false },
    @{ Cell = 'D32'; Value = '0.08656'; Numeric = $true },
    @{ Cell = 'E32'; Value = '  -1.80%  '; Numeric = $false },
    @{ Cell = 'E33'; Value = '  +1.26%  '; Numeric = $false },
    @{ Cell = 'B34'; Value = 'Filecoin'; Numeric = $false },
    @{ Cell = 'C34'; Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; Numeric = $false },
    @{ Cell = 'D34'; Value = '4.438'; Numeric = $true },
    @{ Cell = 'E34'; Value = '  -1.11%  '; Numeric = $false },
    @{ Cell = 'B35'; Value = 'ImmutableX'; Numeric = $false },
    @{ Cell = 'C35'; Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'; Numeric = $false },
    @{ Cell = 'D35'; Value = '0.7223'; Numeric = $true },
    @{ Cell = 'E35'; Value = '  -3.96%  '; Numeric = $false },
    @{ Cell = 'B36'; Value = 'RenderToken'; Numeric = $false },
    @{ Cell = 'C36'; Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'; Numeric = $false },
    @{ Cell = 'D36'; Value = '2.710'; Numeric = $true },
    @{ Cell = 'E36'; Value = '  +5.44%  '; Numeric = $false },
    @{ Cell = 'D37'; Value = '1.103'; Numeric = $true },
    @{ Cell = 'E37'; Value = '  -3.23%  '; Numeric = $false },
    @{ Cell = 'D38'; Value = '1.002'; Numeric = $true },
    @{ Cell = 'E38'; Value = '  -0.17%  '; Numeric = $false },
    @{ Cell = 'D39'; Value = '1.064'; Numeric = $true },
    @{ Cell = 'E39'; Value = '  -2.48%  '; Numeric = $false },
    @{ Cell = 'D40'; Value = '0.01926'; Numeric = $true },
    @{ Cell = 'E40'; Value = '  -1.20%  '; Numeric = $false },
    @{ Cell = 'E41'; Value = '  -1.43%  '; Numeric = $false },
    @{ Cell = 'D42'; Value = '2.853'; Numeric = $true },
    @{ Cell = 'E42'; Value = '  -1.40%  '; Numeric = $false },
    @{ Cell = 'D43'; Value = '0.5061'; Numeric = $true },
    @{ Cell = 'E43'; Value = '  +1.37%  '; Numeric = $false },
    @{ Cell = 'D44'; Value = '6.851'; Numeric = $true },
    @{ Cell = 'E44'; Value = '  -1.68%  '; Numeric = $false },
    @{ Cell = 'D45'; Value = '0.1512'; Numeric = $true },
    @{ Cell = 'E45'; Value = '  -5.71%  '; Numeric = $false },
    @{ Cell = 'D46'; Value = '7.980'; Numeric = $true },
    @{ Cell = 'E46'; Value = '  -4.32%  '; Numeric = $false },
    @{ Cell = 'D47'; Value = '1.002'; Numeric = $true },
    @{ Cell = 'E47'; Value = '  -0.58%  '; Numeric = $false },
    @{ Cell = 'D48'; Value = '0.4620'; Numeric = $true },
    @{ Cell = 'E48'; Value = '  -1.58%  '; Numeric = $false },
    @{ Cell = 'B49'; Value = 'EnergySwap'; Numeric = $false },
    @{ Cell = 'C49'; Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; Numeric = $false },
    @{ Cell = 'D49'; Value = '9.889'; Numeric = $true },
    @{ Cell = 'E49'; Value = '  -2.29%  '; Numeric = $false },
    @{ Cell = 'B50'; Value = 'Quant'; Numeric = $false },
    @{ Cell = 'C50'; Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'; Numeric = $false },
    @{ Cell = 'D50'; Value = '101.02'; Numeric = $true },
    @{ Cell = 'E50'; Value = '  -1.38%  '; Numeric = $false },
    @{ Cell = 'D51'; Value = '1.574'; Numeric = $true },
    @{ Cell = 'E51'; Value = '  -2.47%  '; Numeric = $false }
)

foreach ($update in $updates) {
    $cell = $ws.Range($update.Cell)

    if ($update.Numeric) {
        # Force text storage so the numeric-looking string (e.g. "1.004")
        # isn't coerced into a real number, then restore the default style
        # so we don't leave a stray number-format behind on the cell.
        $cell.NumberFormat = '@'
        $cell.Value = $update.Value
        $cell.Style = 'Normal'
    } else {
        $cell.Value = $update.Value
    }
}
